$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" '35.034.62'
Set-TextCell "E2" '  +1.09%  '
Set-TextCell "D3" '1.851.94'
Set-TextCell "E3" '  +2.31%  '
Set-TextCell "E4" '  -0.02%  '
Set-TextCell "D5" '237.03'
Set-TextCell "E5" '  +3.29%  '
Set-TextCell "E6" '  +0.98%  '
Set-TextCell "D8" '42.50'
Set-TextCell "E8" '  +7.15%  '
Set-TextCell "D9" '0.328'
Set-TextCell "E9" '  +2.48%  '
Set-TextCell "E10" '  +2.17%  '
Set-TextCell "E11" '  +0.56%  '
Set-TextCell "D12" '2.121.48'
Set-TextCell "E12" '  +2.37%  '
Set-TextCell "D13" '11.44'
Set-TextCell "E13" '  +1.94%  '
Set-TextCell "D14" '1.857.01'
Set-TextCell "E14" '  +2.54%  '
Set-TextCell "B15" 'Polkadot'
Set-TextCell "C15" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell "D15" '4.80'
Set-TextCell "E15" '  +4.81%  '
Set-TextCell "B16" 'Polygon'
Set-TextCell "C16" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell "D16" '0.678'
Set-TextCell "E16" '  +2.31%  '
Set-TextCell "D17" '35.034.45'
Set-TextCell "E17" '  +1.44%  '
Set-TextCell "D18" '70.29'
Set-TextCell "E18" '  +1.54%  '
Set-TextCell "D19" '0.0₃0796'
Set-TextCell "E19" '  +1.96%  '
Set-TextCell "D20" '240.70'
Set-TextCell "E20" '  +0.62%  '
Set-TextCell "D21" '12.20'
Set-TextCell "E21" '  +2.93%  '
Set-TextCell "E22" '  +2.92%  '
Set-TextCell "E23" '  +0.02%  '
Set-TextCell "E24" '  +0.98%  '
Set-TextCell "D25" '170.74'
Set-TextCell "E25" '  -1.75%  '
Set-TextCell "D26" '1.89'
Set-TextCell "E26" '  +25.90%  '
Set-TextCell "D27" '7.99'
Set-TextCell "E27" '  +3.60%  '
Set-TextCell "D28" '17.69'
Set-TextCell "E28" '  +2.58%  '
Set-TextCell "E29" '  +0.67%  '
Set-TextCell "D30" '0.0558'
Set-TextCell "E30" '  +2.57%  '
Set-TextCell "E31" '  +0.03%  '
Set-TextCell "E32" '  +0.78%  '
Set-TextCell "E33" '  +3.94%  '
Set-TextCell "E34" '  +23.76%  '
Set-TextCell "E35" '  +12.57%  '
Set-TextCell "D36" '1.34'
Set-TextCell "E36" '  +9.24%  '
Set-TextCell "D37" '0.784'
Set-TextCell "E37" '  +13.93%  '
Set-TextCell "E38" '  +11.70%  '
Set-TextCell "D39" '0.0203'
Set-TextCell "E39" '  +6.35%  '
Set-TextCell "D40" '90.71'
Set-TextCell "E40" '  -0.24%  '
Set-TextCell "D41" '1.350.49'
Set-TextCell "E41" '  +1.39%  '
Set-TextCell "B42" 'InjectiveProtocol'
Set-TextCell "C42" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell "D42" '14.76'
Set-TextCell "E42" '  +3.13%  '
Set-TextCell "B43" 'RenderToken'
Set-TextCell "C43" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell "D43" '2.34'
Set-TextCell "E43" '  +4.45%  '
Set-TextCell "D44" '12.73'
Set-TextCell "E44" '  +53.81%  '
Set-TextCell "E45" '  -0.23%  '
Set-TextCell "D46" '0.0554'
Set-TextCell "E46" '  +7.08%  '
Set-TextCell "E47" '  -0.05%  '
Set-TextCell "E48" '  +6.77%  '
Set-TextCell "D49" '2.034.08'
Set-TextCell "E49" '  +2.12%  '
Set-TextCell "E50" '  +2.59%  '
Set-TextCell "E51" '  +13.44%  '
